# Updates the crypto price/volume table (columns B-E, rows 2-51) on Sheet1
# to reflect the refreshed values from the latest GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value to a cell while keeping it stored as text, even when
# the text looks like a number (e.g. "1.000" or "0.7190"), so Excel does not
# silently coerce it to a numeric value and drop formatting/precision.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "29.483.81"
$ws.Range("E2").Value = "  +0.84%  "

# Row 3
$ws.Range("D3").Value = "1.880.47"

# Row 4
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
Set-TextValue $ws.Range("D5") "0.7190"
$ws.Range("E5").Value = "  +1.90%  "

# Row 6
Set-TextValue $ws.Range("D6") "242.68"
$ws.Range("E6").Value = "  +2.08%  "

# Row 7
Set-TextValue $ws.Range("D7") "1.000"

# Row 8
Set-TextValue $ws.Range("D8") "0.07873"
$ws.Range("E8").Value = "  -1.65%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.3129"
$ws.Range("E9").Value = "  +3.71%  "

# Row 10
Set-TextValue $ws.Range("D10") "25.22"
$ws.Range("E10").Value = "  +7.55%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.08263"
$ws.Range("E11").Value = "  +1.07%  "

# Row 12
$ws.Range("D12").Value = "1.888.85"
$ws.Range("E12").Value = "  +1.10%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.7330"
$ws.Range("E13").Value = "  +4.21%  "

# Row 14
Set-TextValue $ws.Range("D14") "5.295"
$ws.Range("E14").Value = "  +2.21%  "

# Row 15
Set-TextValue $ws.Range("D15") "91.19"
$ws.Range("E15").Value = "  +1.83%  "

# Row 16
$ws.Range("D16").Value = "29.561.63"
$ws.Range("E16").Value = "  +1.17%  "

# Row 17
Set-TextValue $ws.Range("D17") "5.952"
$ws.Range("E17").Value = "  +2.80%  "

# Row 18
Set-TextValue $ws.Range("D18") "247.50"
$ws.Range("E18").Value = "  +4.06%  "

# Row 19
$ws.Range("E19").Value = "  -0.35%  "

# Row 20
Set-TextValue $ws.Range("D20") "13.33"
$ws.Range("E20").Value = "  +0.97%  "

# Row 21
Set-TextValue $ws.Range("D21") "1.000"
$ws.Range("E21").Value = "  -0.02%  "

# Row 22
Set-TextValue $ws.Range("D22") "8.024"
$ws.Range("E22").Value = "  +7.58%  "

# Row 23
Set-TextValue $ws.Range("D23") "1.000"
$ws.Range("E23").Value = "  -0.11%  "

# Row 24
Set-TextValue $ws.Range("D24") "0.1585"
$ws.Range("E24").Value = "  +10.91%  "

# Row 25
Set-TextValue $ws.Range("D25") "164.02"
$ws.Range("E25").Value = "  +0.69%  "

# Row 26
Set-TextValue $ws.Range("D26") "9.055"
$ws.Range("E26").Value = "  +1.85%  "

# Row 27
Set-TextValue $ws.Range("D27") "18.35"
$ws.Range("E27").Value = "  +1.50%  "

# Row 28
Set-TextValue $ws.Range("D28") "1.362"
$ws.Range("E28").Value = "  -4.62%  "

# Row 29
Set-TextValue $ws.Range("D29") "1.495"
$ws.Range("E29").Value = "  +1.35%  "

# Row 30
Set-TextValue $ws.Range("D30") "4.382"
$ws.Range("E30").Value = "  +0.60%  "

# Row 31
Set-TextValue $ws.Range("D31") "4.148"
$ws.Range("E31").Value = "  +3.36%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.05305"
$ws.Range("E32").Value = "  +2.31%  "

# Row 33
Set-TextValue $ws.Range("D33") "1.937"
$ws.Range("E33").Value = "  +1.11%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.205"
$ws.Range("E34").Value = "  +4.17%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.7250"
$ws.Range("E35").Value = "  +1.73%  "

# Row 36
Set-TextValue $ws.Range("D36") "2.676"
$ws.Range("E36").Value = "  +1.07%  "

# Row 37
$ws.Range("E37").Value = "  +1.28%  "

# Row 38
$ws.Range("D38").Value = "1.262.03"
$ws.Range("E38").Value = "  +11.48%  "

# Row 39
Set-TextValue $ws.Range("D39") "2.730"
$ws.Range("E39").Value = "  +0.23%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.9090"
$ws.Range("E40").Value = "  -3.55%  "

# Row 41
Set-TextValue $ws.Range("D41") "74.10"
$ws.Range("E41").Value = "  +5.56%  "

# Row 42
Set-TextValue $ws.Range("D42") "6.109"
$ws.Range("E42").Value = "  +2.90%  "

# Row 43
$ws.Range("E43").Value = "  -0.02%  "

# Row 44
Set-TextValue $ws.Range("D44") "103.83"
$ws.Range("E44").Value = "  +0.95%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.5327"
$ws.Range("E45").Value = "  +0.33%  "

# Row 46  # row 46/47: SynthetixNetwork and RenderToken swapped places
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D46") "1.773"
$ws.Range("E46").Value = "  +0.79%  "

# Row 47  # row 46/47: SynthetixNetwork and RenderToken swapped places
$ws.Range("B47").Value = "SynthetixNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue $ws.Range("D47") "2.931"
$ws.Range("E47").Value = "  +13.56%  "

# Row 48
$ws.Range("E48").Value = "  +2.57%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.4342"
$ws.Range("E49").Value = "  +2.18%  "

# Row 50
Set-TextValue $ws.Range("D50") "9.261"
$ws.Range("E50").Value = "  +1.04%  "

# Row 51
Set-TextValue $ws.Range("D51") "7.080"
$ws.Range("E51").Value = "  +2.14%  "
